$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NVDA row (row 2) values
$ws.Range("B2").Value = 0.004791757726395753
$ws.Range("C2").Value = 0.0007813039001524819

# Update AMAT row (row 3) values
$ws.Range("B3").Value = 0.002288414744028531
$ws.Range("C3").Value = 0.0003758874907965164

# Rows 4 (IDXX) and 5 (CPRT) remain unchanged

# Row 6 becomes CTAS with new values (was TXN)
$ws.Range("A6").Value = "CTAS"
$ws.Range("B6").Value = 0.001280820366687428
$ws.Range("C6").Value = 0.0001526327418421319

# Delete old rows 7-10 (CTAS, SNPS, ADP, PEP) entirely
$ws.Range("A7:C10").EntireRow.Delete()
